$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new row above row 35 (old row 35 "Report FEC" shifts down to 36)
# ---------------------------------------------------------------------------
$ws.Rows(35).EntireRow.Insert()

# ---------------------------------------------------------------------------
# Populate the newly inserted row 35 ("Both TOL inc" / fec scenario row)
# ---------------------------------------------------------------------------
$ws.Range("A35").Formula = "=ROW(A35)-5"
$ws.Range("B35").Value = $false
$ws.Range("C35").Value = $false
$ws.Range("D35").Value = "Both TOL inc"
$ws.Range("E35").Formula = "=E`$5"
$ws.Range("F35").Formula = "=F`$5"
$ws.Range("G35").Formula = "=G`$5"
$ws.Range("H35").Formula = "=H`$5"
$ws.Range("I35").Value = $true
$ws.Range("J35").Formula = "=J`$5"
$ws.Range("K35").Formula = "=K`$5"
$ws.Range("L35").Formula = "=L`$5"
$ws.Range("M35").Formula = "=M`$5"
$ws.Range("N35").Formula = "=N`$5"
$ws.Range("O35").Formula = "=O`$5"
$ws.Range("P35").Formula = "=P`$5"
$ws.Range("Q35").Formula = "=Q`$5"
$ws.Range("R35").Formula = "=R`$5"
$ws.Range("S35").Formula = "=S`$5"
$ws.Range("T35").Formula = "=T`$5"
$ws.Range("U35").Formula = "=U`$5"
$ws.Range("V35").Formula = "=V`$5"
$ws.Range("W35").Formula = "=W`$5"
$ws.Range("X35").Formula = "=X`$5"
$ws.Range("Y35").Formula = "=Y`$5"
$ws.Range("Z35").Formula = "=Z`$5"
$ws.Range("AA35").Formula = "=AA`$5"
$ws.Range("AB35").Formula = "=AB`$5"
$ws.Range("AC35").Formula = "=AC`$5"
$ws.Range("AD35").Formula = "=AD`$5"
$ws.Range("AE35").Formula = "=AE`$5"

# Column D carries a distinct (shaded) style throughout the table - copy it
# down from the row above so the new row matches the rest of the column.
$ws.Range("D34").Copy()
$ws.Range("D35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# View / selection tweaks that came with the edit
# ---------------------------------------------------------------------------
$ws.Range("E20").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("E20").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B37").Select()

Write-Output "done"
